$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.048.39"
$ws.Range("D3").Value = "2.213.85"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "289.14"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "87.90"
$ws.Range("E6").Value = "  +4.42%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.471"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "30.69"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "47.87"
$ws.Range("E12").Value = "  +3.96%  "
$ws.Range("E13").Value = "  +2.37%  "
$ws.Range("D14").Value = "6.46"
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("D15").Value = "2.556.06"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "14.02"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "2.197.86"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "0.728"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "39.966.94"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "11.94"
$ws.Range("E20").Value = "  +14.35%  "
$ws.Range("D21").Value = "0.0₃0886"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "5.81"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").Value = "65.67"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "235.59"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").Value = "1.85"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "22.64"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +5.12%  "
$ws.Range("D30").Value = "9.23"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "152.79"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "32.21"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "4.96"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("D35").Value = "0.0718"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("D38").Value = "15.99"
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("D39").Value = "0.111"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").Value = "0.0997"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("D42").Value = "2.087.54"
$ws.Range("E42").Value = "  +8.18%  "
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "9.87"
$ws.Range("E46").Value = "  +7.05%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "17.52"
$ws.Range("E47").Value = "  +7.24%  "
$ws.Range("D48").Value = "2.67"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("D49").Value = "2.427.54"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "69.65"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "88.72"
$ws.Range("E51").Value = "  -0.08%  "
